$wb = $excel.ActiveWorkbook

# --- metadata sheet: fill in the indicator-specific values in column B ---
$ws = $wb.Worksheets.Item("metadata")

$ws.Range("B2").Value  = "NO_URAQ_001-103"
$ws.Range("B3").Value  = "Air Quality"
$ws.Range("B4").Value  = "Norway"
$ws.Range("B5").Value  = "Europe"
$ws.Range("B6").Value  = "A2 - Chemical State characteristics"
$ws.Range("B7").Value  = "Terrestrial (T)"
$ws.Range("B8").Value  = "T7 Intensive land-use biome"
$ws.Range("B15").Value = "Clappe, S., Czúcz, B."
$ws.Range("B17").Value = "No"
$ws.Range("B18").Value = "TBA"

# --- switch the active sheet / selection to the metadata sheet ---
$ws.Activate() | Out-Null
$ws.Range("B21").Select() | Out-Null
